$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Total_Due (D) and Last_Payment_Date (E) values per the backup refresh.
$ws.Range("D2").Value = 3080.023000000001
$ws.Range("E2").Value = "2025-03-28 18:18:25"

$ws.Range("D3").Value = 74945.20688000001
$ws.Range("E3").Value = "2025-03-28 19:22:09"

$ws.Range("D5").Value = 14700.147
